# Apply the edits described in the commit:
#  1. "[default is no] (yes or no)"  -> "[default is yes] (yes or no)"   (cell H6)
#  2. "use_uncorrected_bedding"      -> "correct_bedding_using_local_dec" (cell H7)
#  3. Active sheet selection moves from G13 to H7
#  4. H7 cell style changes (fill stays the same, alignment changes from
#     "left" to "general" - i.e. style index 8 -> 9 in the original file)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "[default is no]..." text in H6
$ws.Range("H6").Value = "[default is yes] (yes or no)"

# 2. Rename the use_uncorrected_bedding header in H7
$ws.Range("H7").Value = "correct_bedding_using_local_dec"

# 3. Move the active selection to H7
$ws.Range("H7").Select()

# 4. Change H7's horizontal alignment to general (matches style index 9)
$ws.Range("H7").HorizontalAlignment = 1
